# Improved main infographic in chapter 3
# Update the "Methodological decisions ..." caption under the
# "Thesis' integrated experimental strategy" heading so the second
# sentence mentions "research questions and objectives" instead of
# just "research objectives".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$oldText = " underlying the research objectives presented in each main analysis chapter"
$newText = " underlying the research questions and objectives presented in each main analysis chapter"

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)

    if ($sh.HasTextFrame -and $sh.Name -eq "TextBox 65") {
        $tr = $sh.TextFrame.TextRange

        for ($pIdx = 1; $pIdx -le $tr.Paragraphs().Count; $pIdx++) {
            $para = $tr.Paragraphs($pIdx, 1)
            $runs = $para.Runs()

            foreach ($r in $runs) {
                if ($r.Text -eq $oldText) {
                    $r.Text = $newText
                }
            }
        }
    }
}
